$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "ValidLogin"
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"
$ws.Columns("B").ColumnWidth = 10.35
[void]$ws.Range("A3").Select()
$excel.ActiveWindow.Zoom = 205
